$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.127.37"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.998.24"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "245.83"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("D7").Value = "59.85"
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").Value = "0.0805"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "15.00"
$ws.Range("E12").Value = "  +5.25%  "
$ws.Range("D13").Value = "22.50"
$ws.Range("E13").Value = "  +5.46%  "
$ws.Range("D14").Value = "2.292.15"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").Value = "2.008.77"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "37.084.92"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "70.15"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "5.17"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "230.03"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "2.46"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "9.41"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "163.66"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("D29").Value = "19.60"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +12.23%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "0.0654"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("D35").Value = "2.39"
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  -7.09%  "
$ws.Range("D39").Value = "5.36"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "0.0981"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "2.93"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("E44").Value = "  +5.70%  "
$ws.Range("D45").Value = "90.90"
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("D46").Value = "1.366.25"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "7.41"
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("E49").Value = "  +12.20%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "46.25"
$ws.Range("E51").Value = "  +5.19%  "
